{"js": "// Replace the date line and the 25 three-digit-by-one-digit multiplication\n// equations in the table with the values from the \"next day\" worksheet.\n// Every old value in this document is unique, so a plain exact-text\n// search-and-replace (case sensitive, whole match) is sufficient and safe.\nconst replacements = [\n  [\"2024-10-09 Wednesday\", \"2024-10-10 Thursday\"],\n  [\"256\u00d72=512\", \"586\u00d74=2344\"],\n  [\"538\u00d79=4842\", \"486\u00d74=1944\"],\n  [\"867\u00d76=5202\", \"588\u00d76=3528\"],\n  [\"316\u00d72=632\", \"819\u00d76=4914\"],\n  [\"215\u00d76=1290\", \"548\u00d78=4384\"],\n  [\"572\u00d78=4576\", \"767\u00d78=6136\"],\n  [\"607\u00d76=3642\", \"538\u00d78=4304\"],\n  [\"737\u00d75=3685\", \"934\u00d74=3736\"],\n  [\"564\u00d73=1692\", \"212\u00d76=1272\"],\n  [\"877\u00d75=4385\", \"821\u00d75=4105\"],\n  [\"443\u00d76=2658\", \"632\u00d79=5688\"],\n  [\"911\u00d76=5466\", \"749\u00d76=4494\"],\n  [\"847\u00d72=1694\", \"316\u00d77=2212\"],\n  [\"469\u00d77=3283\", \"915\u00d79=8235\"],\n  [\"953\u00d77=6671\", \"149\u00d77=1043\"],\n  [\"356\u00d77=2492\", \"599\u00d75=2995\"],\n  [\"622\u00d78=4976\", \"282\u00d76=1692\"],\n  [\"189\u00d79=1701\", \"812\u00d73=2436\"],\n  [\"619\u00d72=1238\", \"534\u00d78=4272\"],\n  [\"909\u00d78=7272\", \"431\u00d76=2586\"],\n  [\"745\u00d78=5960\", \"469\u00d79=4221\"],\n  [\"143\u00d75=715\", \"680\u00d78=5440\"],\n  [\"784\u00d79=7056\", \"142\u00d78=1136\"],\n  [\"565\u00d79=5085\", \"132\u00d77=924\"],\n  [\"838\u00d78=6704\", \"587\u00d72=1174\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const hits = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  hits.load(\"items\");\n  await context.sync();\n\n  for (const hit of hits.items) {\n    hit.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 three-digit-by-one-digit multiplication\n# equations in the table with the values from the \"next day\" worksheet.\n# Every old value in this document is unique, so a plain exact-text\n# Find/Replace-all (case sensitive, whole document) is sufficient and safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-10-09 Wednesday\", \"2024-10-10 Thursday\"),\n    @(\"256\u00d72=512\", \"586\u00d74=2344\"),\n    @(\"538\u00d79=4842\", \"486\u00d74=1944\"),\n    @(\"867\u00d76=5202\", \"588\u00d76=3528\"),\n    @(\"316\u00d72=632\", \"819\u00d76=4914\"),\n    @(\"215\u00d76=1290\", \"548\u00d78=4384\"),\n    @(\"572\u00d78=4576\", \"767\u00d78=6136\"),\n    @(\"607\u00d76=3642\", \"538\u00d78=4304\"),\n    @(\"737\u00d75=3685\", \"934\u00d74=3736\"),\n    @(\"564\u00d73=1692\", \"212\u00d76=1272\"),\n    @(\"877\u00d75=4385\", \"821\u00d75=4105\"),\n    @(\"443\u00d76=2658\", \"632\u00d79=5688\"),\n    @(\"911\u00d76=5466\", \"749\u00d76=4494\"),\n    @(\"847\u00d72=1694\", \"316\u00d77=2212\"),\n    @(\"469\u00d77=3283\", \"915\u00d79=8235\"),\n    @(\"953\u00d77=6671\", \"149\u00d77=1043\"),\n    @(\"356\u00d77=2492\", \"599\u00d75=2995\"),\n    @(\"622\u00d78=4976\", \"282\u00d76=1692\"),\n    @(\"189\u00d79=1701\", \"812\u00d73=2436\"),\n    @(\"619\u00d72=1238\", \"534\u00d78=4272\"),\n    @(\"909\u00d78=7272\", \"431\u00d76=2586\"),\n    @(\"745\u00d78=5960\", \"469\u00d79=4221\"),\n    @(\"143\u00d75=715\", \"680\u00d78=5440\"),\n    @(\"784\u00d79=7056\", \"142\u00d78=1136\"),\n    @(\"565\u00d79=5085\", \"132\u00d77=924\"),\n    @(\"838\u00d78=6704\", \"587\u00d72=1174\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
